# Add a new trilateration data point (row 13): an agent position plus the
# three anchors/distances used for the linear-algebra intersection calc,
# mirroring rows 2-12's layout (idx, agent_x, agent_y, anchor1_x/y, dist1,
# anchor2_x/y, dist2, anchor3_x/y, dist3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 1
$ws.Range("F13").Formula = "=SQRT((D13-`$B13)^2+(E13-`$C13)^2)"
$ws.Range("G13").Value = 13
$ws.Range("H13").Value = 12
$ws.Range("I13").Formula = "=SQRT((G13-`$B13)^2+(H13-`$C13)^2)"
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 12
$ws.Range("L13").Formula = "=SQRT((J13-`$B13)^2+(K13-`$C13)^2)"
